$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.223.60'
$ws.Range("E2").Value = '  -1.07%  '
$ws.Range("D3").Value = '3.529.79'
$ws.Range("E3").Value = '  +0.48%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = '''607.09'
$ws.Range("E5").Value = '  -0.24%  '
$ws.Range("D6").Value = '''143.19'
$ws.Range("E6").Value = '  -3.09%  '
$ws.Range("D7").Value = '3.526.85'
$ws.Range("E7").Value = '  +0.41%  '
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("E9").Value = '  +0.33%  '
$ws.Range("D10").Value = '''0.136'
$ws.Range("E10").Value = '  -4.33%  '
$ws.Range("D11").Value = '''8.05'
$ws.Range("E11").Value = '  +1.09%  '
$ws.Range("E12").Value = '  -2.88%  '
$ws.Range("D13").Value = '4.128.99'
$ws.Range("E13").Value = '  +0.52%  '
$ws.Range("E14").Value = '  -4.78%  '
$ws.Range("D15").Value = '''30.00'
$ws.Range("E15").Value = '  -5.87%  '
$ws.Range("D16").Value = '3.531.71'
$ws.Range("E16").Value = '  +0.56%  '
$ws.Range("D17").Value = '66.316.01'
$ws.Range("E17").Value = '  -1.01%  '
$ws.Range("E18").Value = '  -0.61%  '
$ws.Range("E19").Value = '  +2.39%  '
$ws.Range("D20").Value = '''6.20'
$ws.Range("E20").Value = '  -3.98%  '
$ws.Range("D21").Value = '''14.92'
$ws.Range("E21").Value = '  -2.51%  '
$ws.Range("D22").Value = '''425.30'
$ws.Range("E22").Value = '  -2.91%  '
$ws.Range("D23").Value = '''0.600'
$ws.Range("E23").Value = '  -1.30%  '
$ws.Range("D24").Value = '''78.72'
$ws.Range("E24").Value = '  -0.76%  '
$ws.Range("D25").Value = '3.673.55'
$ws.Range("E25").Value = '  +0.38%  '
$ws.Range("E26").Value = '  +0.01%  '
$ws.Range("E27").Value = '  -1.08%  '
$ws.Range("D28").Value = '''8.05'
$ws.Range("E28").Value = '  -2.57%  '
$ws.Range("D29").Value = '''9.14'
$ws.Range("E29").Value = '  -6.38%  '
$ws.Range("E30").Value = '  -1.59%  '
$ws.Range("D31").Value = '''1.00'
$ws.Range("E31").Value = '  +0.11%  '
$ws.Range("D32").Value = '''0.161'
$ws.Range("E32").Value = '  -4.27%  '
$ws.Range("E33").Value = '  -6.28%  '
$ws.Range("D34").Value = '''25.29'
$ws.Range("E34").Value = '  -0.75%  '
$ws.Range("D35").Value = '3.521.50'
$ws.Range("E36").Value = '  -0.04%  '
$ws.Range("D37").Value = '''1.75'
$ws.Range("E37").Value = '  -3.06%  '
$ws.Range("D38").Value = '''7.81'
$ws.Range("E38").Value = '  -2.60%  '
$ws.Range("E39").Value = '  -5.74%  '
$ws.Range("E40").Value = '  +0.03%  '
$ws.Range("D41").Value = '''171.38'
$ws.Range("E41").Value = '  -0.95%  '
$ws.Range("D42").Value = '''0.0856'
$ws.Range("E42").Value = '  -4.02%  '
$ws.Range("D43").Value = '''5.18'
$ws.Range("E43").Value = '  -4.43%  '
$ws.Range("D44").Value = '''0.891'
$ws.Range("E44").Value = '  -0.45%  '
$ws.Range("E45").Value = '  -8.99%  '
$ws.Range("D46").Value = '''45.38'
$ws.Range("E46").Value = '  -1.65%  '
$ws.Range("D47").Value = '''25.96'
$ws.Range("E47").Value = '  -6.42%  '
$ws.Range("E48").Value = '  -5.23%  '
$ws.Range("E49").Value = '  -2.39%  '
$ws.Range("D50").Value = '''7.15'
$ws.Range("E50").Value = '  -4.17%  '
$ws.Range("D51").Value = '''0.944'
$ws.Range("E51").Value = '  -4.68%  '
